# atualização do modulo recondutoramento
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the remaining data row with the new value
$ws.Range("A2").Value = 2000056127

# Remove the now-obsolete rows 3 and 4 (2000033006 / 2000062283)
$ws.Rows("3:4").Delete()

# Reflect the new selection/active cell used when the file was last saved
$ws.Range("C7").Select()
